$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 666.8570999999999
$ws.Range("I41").Value = 850.25
$ws.Range("J41").Value = 593.5
$ws.Range("K41").Value = 850.25
$ws.Range("L41").Value = 593.5
$ws.Range("M41").Value = -410.25
$ws.Range("N41").Value = -1473.5
$ws.Range("H62").Value = 3495
$ws.Range("H65").Value = 3495
$ws.Range("H112").Value = 1737.5217
$ws.Range("I112").Value = 1166.6666
$ws.Range("J112").Value = 1823.15
$ws.Range("K112").Value = 3499.9998
$ws.Range("L112").Value = 5469.450000000001
$ws.Range("M112").Value = -2391.9998
$ws.Range("N112").Value = -7685.450000000001
$ws.Range("H135").Value = 1103.0667
$ws.Range("I135").Value = 637.7143
$ws.Range("J135").Value = 1510.25
$ws.Range("K135").Value = 5739.428699999999
$ws.Range("L135").Value = 13592.25
$ws.Range("M135").Value = -3204.428699999999
$ws.Range("N135").Value = -18662.25
$ws.Range("H137").Value = 1417.3334
$ws.Range("J137").Value = 1560
$ws.Range("L137").Value = 4680
$ws.Range("N137").Value = -9780
$ws.Range("H138").Value = 3376.8667
$ws.Range("J138").Value = 4286.125
$ws.Range("L138").Value = 12858.375
$ws.Range("N138").Value = -23138.375
$ws.Range("H140").Value = 64497.5
$ws.Range("J140").Value = 64497.5
$ws.Range("L140").Value = 64497.5
$ws.Range("N140").Value = -74857.5
$ws.Range("H141").Value = 4428.4287
$ws.Range("I141").Value = 4239.8
$ws.Range("J141").Value = 4900
$ws.Range("K141").Value = 12719.4
$ws.Range("L141").Value = 14700
$ws.Range("M141").Value = -7539.400000000001
$ws.Range("N141").Value = -25060

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()
$ws.Range("H61").Value = 2058.2942
$ws.Range("I61").Value = 1830.2307
$ws.Range("J61").Value = 2799.5
$ws.Range("K61").Value = 1830.2307
$ws.Range("L61").Value = 2799.5
$ws.Range("M61").Value = -1618.2307
$ws.Range("N61").Value = -3223.5
$ws.Range("H136").Value = 2058.2942
$ws.Range("I136").Value = 1830.2307
$ws.Range("J136").Value = 2799.5
$ws.Range("K136").Value = 5490.6921
$ws.Range("L136").Value = 8398.5
$ws.Range("M136").Value = -2940.6921
$ws.Range("N136").Value = -13498.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 60555.555
$ws.Range("J132").Value = 60555.555
$ws.Range("L132").Value = 60555.555
$ws.Range("N132").Value = -70675.55499999999
$ws.Range("H141").Value = 30000
$ws.Range("I141").Value = 30000
$ws.Range("K141").Value = 30000
$ws.Range("M141").Value = -24820

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 3978.6667
$ws.Range("I86").Value = 3425
$ws.Range("J86").Value = 4421.6
$ws.Range("K86").Value = 3425
$ws.Range("L86").Value = 4421.6
$ws.Range("M86").Value = -2302
$ws.Range("N86").Value = -6667.6
$ws.Range("H89").Value = 3978.6667
$ws.Range("I89").Value = 3425
$ws.Range("J89").Value = 4421.6
$ws.Range("K89").Value = 17125
$ws.Range("L89").Value = 22108
$ws.Range("M89").Value = -11509
$ws.Range("N89").Value = -33340

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1242.6061
$ws.Range("I5").Value = 1107.7333
$ws.Range("J5").Value = 1355
$ws.Range("K5").Value = 3323.199900000001
$ws.Range("L5").Value = 4065
$ws.Range("M5").Value = -3211.199900000001
$ws.Range("N5").Value = -4289
$ws.Range("H37").Value = 647442.7
$ws.Range("J37").Value = 647442.7
$ws.Range("L37").Value = 1942328.1
$ws.Range("N37").Value = -1942552.1
$ws.Range("H107").Value = 515378
$ws.Range("I107").Value = 550.5714
$ws.Range("K107").Value = 1651.7142
$ws.Range("M107").Value = 268.2857999999999
$ws.Range("H131").Value = 788.95
$ws.Range("I131").Value = 455.75
$ws.Range("J131").Value = 802.8333
$ws.Range("K131").Value = 1367.25
$ws.Range("L131").Value = 2408.4999
$ws.Range("M131").Value = 3672.75
$ws.Range("N131").Value = -12488.4999
$ws.Range("H135").Value = 1242.6061
$ws.Range("I135").Value = 1107.7333
$ws.Range("J135").Value = 1355
$ws.Range("K135").Value = 9969.599700000001
$ws.Range("L135").Value = 12195
$ws.Range("M135").Value = -7434.599700000001
$ws.Range("N135").Value = -17265
$ws.Range("H137").Value = 15927953
$ws.Range("I137").Value = 127226.25
$ws.Range("J137").Value = 25651478
$ws.Range("K137").Value = 381678.75
$ws.Range("L137").Value = 76954434
$ws.Range("M137").Value = -376578.75
$ws.Range("N137").Value = -76964634
$ws.Range("H138").Value = 13211
$ws.Range("I138").Value = 18316.666
$ws.Range("J138").Value = 2999.6667
$ws.Range("K138").Value = 54949.99800000001
$ws.Range("L138").Value = 8999.000100000001
$ws.Range("M138").Value = -49809.99800000001
$ws.Range("N138").Value = -19279.0001
$ws.Range("H140").Value = 4652.25
$ws.Range("I140").Value = 5958.8096
$ws.Range("J140").Value = 2157.9092
$ws.Range("K140").Value = 17876.4288
$ws.Range("L140").Value = 6473.7276
$ws.Range("M140").Value = -12696.4288
$ws.Range("N140").Value = -16833.7276
$ws.Range("H141").Value = 7918.125
$ws.Range("I141").Value = 8603.076999999999
$ws.Range("J141").Value = 4950
$ws.Range("K141").Value = 25809.231
$ws.Range("L141").Value = 14850
$ws.Range("M141").Value = -20629.231
$ws.Range("N141").Value = -25210

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4048.7144
$ws.Range("I102").Value = 3178.4
$ws.Range("J102").Value = 4532.222
$ws.Range("K102").Value = 3178.4
$ws.Range("L102").Value = 4532.222
$ws.Range("M102").Value = -1556.4
$ws.Range("N102").Value = -7776.222
$ws.Range("H122").Value = 745.3889
$ws.Range("I122").Value = 744.8125
$ws.Range("J122").Value = 750
$ws.Range("K122").Value = 2234.4375
$ws.Range("L122").Value = 2250
$ws.Range("M122").Value = 215.5625
$ws.Range("N122").Value = -7150
$ws.Range("H140").Value = 105937.375
$ws.Range("J140").Value = 105937.375
$ws.Range("L140").Value = 105937.375
$ws.Range("N140").Value = -116297.375
$ws.Range("H141").Value = 47494.75
$ws.Range("J141").Value = 47494.75
$ws.Range("L141").Value = 47494.75
$ws.Range("N141").Value = -57854.75

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1709.1305
$ws.Range("J61").Value = 1900.4546
$ws.Range("L61").Value = 1900.4546
$ws.Range("N61").Value = -2304.4546
$ws.Range("H113").Value = 1709.1305
$ws.Range("J113").Value = 1900.4546
$ws.Range("L113").Value = 1900.4546
$ws.Range("N113").Value = -6240.4546

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 842.5
$ws.Range("I113").Value = 697.375
$ws.Range("J113").Value = 987.625
$ws.Range("K113").Value = 2092.125
$ws.Range("L113").Value = 2962.875
$ws.Range("M113").Value = 77.875
$ws.Range("N113").Value = -7302.875
$ws.Range("H136").Value = 2150.5454
$ws.Range("I136").Value = 631.38464
$ws.Range("J136").Value = 4344.8887
$ws.Range("K136").Value = 1894.15392
$ws.Range("L136").Value = 13034.6661
$ws.Range("M136").Value = 655.84608
$ws.Range("N136").Value = -18134.6661
